# Auto-generated script applying 2023-04-27 weekly crime data update
# to the J column (2023 running total) across Citywide Totals, By Neighborhood,
# and all affected per-neighborhood worksheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 2142
$ws.Range("J3").Value = 2240
$ws.Range("J4").Value = 507
$ws.Range("J5").Value = 162
$ws.Range("J6").Value = 2834
$ws.Range("J7").Value = 7885

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("J2").Value = 31
$ws.Range("J7").Value = 88

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J2").Value = 78
$ws.Range("J3").Value = 95
$ws.Range("J7").Value = 265

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J3").Value = 12
$ws.Range("J7").Value = 53

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J2").Value = 61
$ws.Range("J3").Value = 54
$ws.Range("J6").Value = 74
$ws.Range("J7").Value = 204

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J2").Value = 60
$ws.Range("J4").Value = 33
$ws.Range("J6").Value = 81
$ws.Range("J7").Value = 235
$ws.Range("J8").Value = 498
$ws.Range("J10").Value = 48
$ws.Range("J11").Value = 107
$ws.Range("J19").Value = 259
$ws.Range("J20").Value = 163
$ws.Range("J23").Value = 72
$ws.Range("J27").Value = 47
$ws.Range("J29").Value = 442
$ws.Range("J31").Value = 53
$ws.Range("J33").Value = 328
$ws.Range("J36").Value = 118
$ws.Range("J37").Value = 265
$ws.Range("J41").Value = 50
$ws.Range("J42").Value = 303
$ws.Range("J43").Value = 77
$ws.Range("J44").Value = 65
$ws.Range("J48").Value = 71
$ws.Range("J50").Value = 44
$ws.Range("J52").Value = 194
$ws.Range("J53").Value = 75
$ws.Range("J56").Value = 9
$ws.Range("J63").Value = 33
$ws.Range("J64").Value = 51
$ws.Range("J65").Value = 204
$ws.Range("J66").Value = 18
$ws.Range("J71").Value = 33
$ws.Range("J72").Value = 30
$ws.Range("J76").Value = 115
$ws.Range("J78").Value = 108
$ws.Range("J83").Value = 188
$ws.Range("J85").Value = 370
$ws.Range("J93").Value = 38
$ws.Range("J96").Value = 88
$ws.Range("J98").Value = 50
$ws.Range("J101").Value = 7885

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J3").Value = 63
$ws.Range("J6").Value = 56
$ws.Range("J7").Value = 188

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J2").Value = 83
$ws.Range("J3").Value = 98
$ws.Range("J4").Value = 16
$ws.Range("J6").Value = 117
$ws.Range("J7").Value = 328

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J3").Value = 30
$ws.Range("J6").Value = 77

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J2").Value = 131
$ws.Range("J3").Value = 153
$ws.Range("J6").Value = 120
$ws.Range("J7").Value = 442

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J3").Value = 70
$ws.Range("J6").Value = 98
$ws.Range("J7").Value = 259

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("J4").Value = 6
$ws.Range("J7").Value = 65

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("J2").Value = 13
$ws.Range("J7").Value = 71

$ws = $wb.Worksheets.Item("River North")
$ws.Range("J2").Value = 14
$ws.Range("J7").Value = 115

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J2").Value = 89
$ws.Range("J3").Value = 146
$ws.Range("J6").Value = 104
$ws.Range("J7").Value = 370

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("J2").Value = 23
$ws.Range("J4").Value = 3
$ws.Range("J7").Value = 81

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("J6").Value = 24
$ws.Range("J7").Value = 50

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J6").Value = 155
$ws.Range("J7").Value = 303

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("J2").Value = 16
$ws.Range("J7").Value = 48

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J4").Value = 15
$ws.Range("J7").Value = 108

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("J4").Value = 6
$ws.Range("J7").Value = 72

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("J2").Value = 13
$ws.Range("J5").Value = 3
$ws.Range("J7").Value = 51

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J6").Value = 44
$ws.Range("J7").Value = 163

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J2").Value = 40
$ws.Range("J7").Value = 118

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("J6").Value = 13
$ws.Range("J7").Value = 38

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J2").Value = 44
$ws.Range("J4").Value = 9
$ws.Range("J6").Value = 81
$ws.Range("J7").Value = 194

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("J4").Value = 5
$ws.Range("J7").Value = 50

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("J3").Value = 13
$ws.Range("J7").Value = 44

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("J3").Value = 4
$ws.Range("J7").Value = 18

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("J6").Value = 41
$ws.Range("J7").Value = 107

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("J3").Value = 19
$ws.Range("J7").Value = 60

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 159
$ws.Range("J6").Value = 144
$ws.Range("J7").Value = 498

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("J3").Value = 12
$ws.Range("J7").Value = 47

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("J6").Value = 48
$ws.Range("J7").Value = 77

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("J2").Value = 12
$ws.Range("J7").Value = 75

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("J6").Value = 18
$ws.Range("J7").Value = 33

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("J2").Value = 10
$ws.Range("J7").Value = 30

$ws = $wb.Worksheets.Item("Magnificent Mile")
$ws.Range("J2").Value = 1
$ws.Range("J7").Value = 9

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J6").Value = 81
$ws.Range("J7").Value = 235

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("J2").Value = 12
$ws.Range("J7").Value = 33
